$d = $word.ActiveDocument

# The empty paragraph right after "Define the problem..." needs to receive
# the answer text, and the "_GoBack" bookmark needs to move from its old
# (now-empty) paragraph into this paragraph, collapsed right after the new
# text (no visible text wrapped by the bookmark).

$answerText = "The goal is for the man to transport all of the items across the river.  The problem is that the boat is not big enough to transport all of the items at the same time."

# Locate the target empty paragraph: the one immediately following the
# paragraph containing "Define the problem".
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Define the problem") {
        $targetPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

# Write the answer text with one extra placeholder character appended. This
# avoids placing a collapsed bookmark at the very last character slot of a
# paragraph (immediately before the paragraph mark), which this runtime
# mishandles. The placeholder is stripped again right after the bookmark is
# anchored.
$targetPara.Range.Text = $answerText + "X"

$targetParaAfter = $targetPara
# Re-fetch range/paragraph to get fresh Start/End after the text write.
$freshRange = $targetParaAfter.Range
$bookmarkPos = $freshRange.End - 2  # right after the real text, before "X"

# Remove any existing "_GoBack" bookmark (its old location becomes a plain
# empty paragraph) before re-adding it at the new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the placeholder character now that the bookmark is anchored.
$placeholderRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$placeholderRange.Delete()
